# Applies the "debug some error when send code of item to electronic board"
# edit: renumbers/re-codes the LowAmper rooms and adds a new item row, fixes
# the HighAmper zone data, and updates the last-used selection on each sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "LowAmper"
# ---------------------------------------------------------------------------
$low = $wb.Worksheets.Item("LowAmper")

# Row 18 - code fixed, item renamed
$low.Cells.Item(18, 1).Value = 18
$low.Cells.Item(18, 2).Value = 1
$low.Cells.Item(18, 3).Value = "TV Room"
$low.Cells.Item(18, 4).Value = 1
$low.Cells.Item(18, 5).Value = "lamp-1"
$low.Cells.Item(18, 6).Value = "Lamp"
$low.Cells.Item(18, 7).Value = "null"

# Row 19 - now the zone's lamp/lustre item
$low.Cells.Item(19, 1).Value = 19
$low.Cells.Item(19, 2).Value = 9
$low.Cells.Item(19, 3).Value = "اتاق 3"
$low.Cells.Item(19, 4).Value = 1
$low.Cells.Item(19, 5).Value = "لوستر"
$low.Cells.Item(19, 6).Value = "Lamp"
$low.Cells.Item(19, 7).Value = "null"

# Row 20 - right curtain
$low.Cells.Item(20, 1).Value = 20
$low.Cells.Item(20, 2).Value = 9
$low.Cells.Item(20, 3).Value = "اتاق 3"
$low.Cells.Item(20, 4).Value = 2
$low.Cells.Item(20, 5).Value = "راست"
$low.Cells.Item(20, 6).Value = "Curtain"
$low.Cells.Item(20, 7).Value = "null"

# Row 21 - left curtain (new row, previously the aircondition row)
$low.Cells.Item(21, 1).Value = 21
$low.Cells.Item(21, 2).Value = 9
$low.Cells.Item(21, 3).Value = "اتاق 3"
$low.Cells.Item(21, 4).Value = 3
$low.Cells.Item(21, 5).Value = "چپ"
$low.Cells.Item(21, 6).Value = "Curtain"
$low.Cells.Item(21, 7).Value = "null"

# Row 22 - new row, aircondition item
$low.Cells.Item(22, 1).Value = 22
$low.Cells.Item(22, 2).Value = 9
$low.Cells.Item(22, 3).Value = "اتاق 3"
$low.Cells.Item(22, 4).Value = 4
$low.Cells.Item(22, 5).Value = "کولر 3"
$low.Cells.Item(22, 6).Value = "Aircondition"
$low.Cells.Item(22, 7).Value = "Slow"

[void]$low.Range("A14:J21").Select()

# ---------------------------------------------------------------------------
# Sheet "HighAmper"
# ---------------------------------------------------------------------------
$high = $wb.Worksheets.Item("HighAmper")

$high.Cells.Item(2, 2).Value = 9
$high.Cells.Item(2, 3).Value = "اتاق 3"
$high.Cells.Item(2, 5).Value = "کولر 3"

$high.Cells.Item(3, 2).Value = 9
$high.Cells.Item(3, 3).Value = "اتاق 3"
$high.Cells.Item(3, 5).Value = "کولر 3"

[void]$high.Range("G19").Select()
